# Weekly price update: a new price observation (dated 45142, i.e. 2023-08-04)
# is inserted for "Espinaca" at "Terminal La Palmera de La Serena", pushing
# every existing record from row 398 down by one row (398->399, ..., 455->456).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 398; Excel shifts rows 398:455 down to 399:456
# and extends the used range to A1:R456 automatically.
$ws.Rows(398).Insert()

# Populate the newly inserted row 398 with the new weekly record.
$ws.Range("A398").Value = 8
$ws.Range("B398").Value = "Terminal La Palmera de La Serena"
$ws.Range("C398").Value = "Coquimbo"
$ws.Range("D398").Value = 45142
$ws.Range("E398").Value = 4
$ws.Range("F398").Value = 100112012
$ws.Range("G398").Value = "Espinaca"
$ws.Range("H398").Value = "Sin especificar"
$ws.Range("I398").Value = "Primera"
$ws.Range("J398").Value = 1600
$ws.Range("K398").Value = 500
$ws.Range("L398").Value = 600
$ws.Range("M398").Value = 550
$ws.Range("N398").Value = "`$/atado 300 a 500 gramos"
$ws.Range("O398").Value = "Provincia del Elquí"
$ws.Range("P398").Value = 1100
$ws.Range("Q398").Value = 0.5
$ws.Range("R398").Value = "Hortaliza"

# Keep the date cell formatted the same way as the other date cells in column D.
$ws.Range("D398").NumberFormat = $ws.Range("D399").NumberFormat
